$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 155.9807553291321
$ws.Range("C2").Value = 58.87547225374421
$ws.Range("D2").Value = 2.226694631576538
$ws.Range("E2").Value = 3.088080238127659
$ws.Range("B3").Value = 417.4680068969727
$ws.Range("C3").Value = 131.4527615176957
$ws.Range("D3").Value = 4.612991762161255
$ws.Range("E3").Value = 2.038105918452289
$ws.Range("B4").Value = 831.4452156543732
$ws.Range("C4").Value = 147.1369575718296
$ws.Range("D4").Value = 7.794174575805664
$ws.Range("E4").Value = 0.697301385337306
$ws.Range("B5").Value = 377.1240797996521
$ws.Range("C5").Value = 4.488915750976273
$ws.Range("D5").Value = 6.424685859680176
$ws.Range("E5").Value = 0.6204430465207618
$ws.Range("B6").Value = 747.3299785614014
$ws.Range("C6").Value = 1.435729779805043
$ws.Range("D6").Value = 7.767447471618652
$ws.Range("E6").Value = 0.5263217839837919
$ws.Range("B7").Value = 1503.241074752808
$ws.Range("C7").Value = 8.216632186597469
$ws.Range("D7").Value = 7.806480121612549
$ws.Range("E7").Value = 0.4952123317853763
$ws.Range("B8").Value = 505.8822907924652
$ws.Range("C8").Value = 2.412524486109187
$ws.Range("D8").Value = 7.771378517150879
$ws.Range("E8").Value = 0.5204026821273751
$ws.Range("B9").Value = 981.9365540504456
$ws.Range("C9").Value = 3.121151172707056
$ws.Range("D9").Value = 7.429174852371216
$ws.Range("E9").Value = 0.5210085033222966
$ws.Range("B10").Value = 1933.971311092377
$ws.Range("C10").Value = 10.17687054758597
$ws.Range("D10").Value = 8.29528088569641
$ws.Range("E10").Value = 0.5526351660262588
$ws.Range("B11").Value = 364.2634794235229
$ws.Range("C11").Value = 10.36124568332492
$ws.Range("D11").Value = 8.166208696365356
$ws.Range("E11").Value = 1.09103588865048
$ws.Range("B12").Value = 738.1327583312989
$ws.Range("C12").Value = 17.185077873209
$ws.Range("D12").Value = 8.48124794960022
$ws.Range("E12").Value = 0.1751275123734161
$ws.Range("B13").Value = 1443.930177211761
$ws.Range("C13").Value = 14.65252708229283
$ws.Range("D13").Value = 7.77112889289856
$ws.Range("E13").Value = 0.734735051367882
$ws.Range("B14").Value = 531.909878540039
$ws.Range("C14").Value = 10.00108584064636
$ws.Range("D14").Value = 8.041028881072998
$ws.Range("E14").Value = 0.4545830724695889
$ws.Range("B15").Value = 1054.598804473877
$ws.Range("C15").Value = 13.04705813615584
$ws.Range("D15").Value = 7.271795701980591
$ws.Range("E15").Value = 1.090471930160971
$ws.Range("B16").Value = 2102.181675291061
$ws.Range("C16").Value = 24.58502204365717
$ws.Range("D16").Value = 7.786889791488647
$ws.Range("E16").Value = 0.891355749180525
$ws.Range("B17").Value = 702.1941849708558
$ws.Range("C17").Value = 9.525333471544933
$ws.Range("D17").Value = 7.351981019973755
$ws.Range("E17").Value = 0.5475775561768333
$ws.Range("B18").Value = 1386.648744726181
$ws.Range("C18").Value = 16.12656439762762
$ws.Range("D18").Value = 7.236619520187378
$ws.Range("E18").Value = 1.124575695069342
$ws.Range("B19").Value = 2747.567058753967
$ws.Range("C19").Value = 38.98479698671913
$ws.Range("D19").Value = 8.358428239822388
$ws.Range("E19").Value = 0.9187743836990696
$ws.Range("B20").Value = 403.0114535331726
$ws.Range("C20").Value = 9.556101852337729
$ws.Range("D20").Value = 7.453378868103028
$ws.Range("E20").Value = 0.6873262598454917
$ws.Range("B21").Value = 790.3326588153839
$ws.Range("C21").Value = 15.77419469972016
$ws.Range("D21").Value = 8.232116937637329
$ws.Range("E21").Value = 0.4802826904648368
$ws.Range("B22").Value = 1556.118911647797
$ws.Range("C22").Value = 22.65195437510166
$ws.Range("D22").Value = 8.390181016921996
$ws.Range("E22").Value = 0.4493497698688949
$ws.Range("B23").Value = 574.2722779750824
$ws.Range("C23").Value = 14.12138459417679
$ws.Range("D23").Value = 8.251580953598022
$ws.Range("E23").Value = 0.0506019854806141
$ws.Range("B24").Value = 1145.251269388199
$ws.Range("C24").Value = 25.60517341490359
$ws.Range("D24").Value = 8.196055459976197
$ws.Range("E24").Value = 0.4691277559583056
$ws.Range("B25").Value = 2281.5448802948
$ws.Range("C25").Value = 50.18006827834968
$ws.Range("D25").Value = 8.834143114089965
$ws.Range("E25").Value = 0.08239113847584993
$ws.Range("B26").Value = 758.4038465499877
$ws.Range("C26").Value = 22.60167229475072
$ws.Range("D26").Value = 8.271647024154664
$ws.Range("E26").Value = 0.2250530510250268
$ws.Range("B27").Value = 1477.770034837723
$ws.Range("C27").Value = 26.73667225521547
$ws.Range("D27").Value = 7.51599702835083
$ws.Range("E27").Value = 1.099449395601465
$ws.Range("B28").Value = 1950.318675088882
$ws.Range("C28").Value = 318.2246469208184
$ws.Range("D28").Value = 2.339125490188599
$ws.Range("E28").Value = 1.576596433059863
